$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010 and 2010-18")

# The commit logs a new CW3M build (C618): insert a fresh row above the
# current last log entry (row 109) so rows 109.. shift down by one and a
# blank spacer row remains before the older "2010-19" comparison block.
$ws.Rows.Item(109).Insert()

# The new entry's metric columns (D:S) start out identical to the most
# recent existing entry (now row 108) - duplicate that row's values down
# into the newly inserted row 109.
$ws.Range("A108:S108").Copy()
$ws.Range("A109:S109").PasteSpecial(-4163)  # xlPasteValues
$ws.Application.CutCopyMode = 0

# Re-label the new row with the C618 run identifiers.
$ws.Cells.Item(109, 1).Value = "CW3M C618"
$ws.Cells.Item(109, 2).Value = "Demo_Baseline_2010-18_C618"

# Leave the selection on the new row, matching where the edit was made.
$ws.Cells.Item(109, 18).Select()

# Record the workbook window's on-screen position/size at save time.
$win = $ws.Application.ActiveWindow
$win.Left = 28680
$win.Top = -7425
$win.Width = 29040
$win.Height = 17640
